$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Population_name" values in column B were simplified across every
# data row: "Test - Test - 10/30/2020" -> "Test - Test".
$ws.Range("B2").Value = "Test - Test"
$ws.Range("B4").Value = "Test - Test"
$ws.Range("B8").Value = "Test - Test"
$ws.Range("B12").Value = "Test - Test"
$ws.Range("B16").Value = "Test - Test"

# Column B's best-fit width shrinks now that its longest entry is shorter.
$ws.Columns.Item(2).ColumnWidth = 14.65

# The last thing selected/saved was cell C16.
$ws.Range("C16").Select()
